$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Antoine Creek-Lower DS" to "Antoine Creek-Lower" in row 19 (EDT AU column, B)
$ws.Range("B19").Value = "Antoine Creek-Lower"

# Add a new row of data (row 43): Lower Omak Creek / Omak Creek-Lower DS / 1 / 1 / 170200061905
$ws.Range("A43").Value = "Lower Omak Creek"
$ws.Range("B43").Value = "Omak Creek-Lower DS"
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 170200061905
